$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header column (H1), copying the style of the existing
# header cell (G1) so it matches the other headers (bold, centered, bordered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the data value for the new column (H2).
$ws.Range("H2").Value = 0
